# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text block with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$texto = $cellA1.Value2
$texto = $texto.Replace("1000 Bs = 10.3 = 42739.55 pesos", "1000 Bs = 10.24 = 42413.69 pesos")
$texto = $texto.Replace("42739.55 pesos = 10.28 = 976.73 Bs", "42413.69 pesos = 10.19 = 956.92 Bs")
$cellA1.Value = $texto

# --- tasas: update the raw rate figures in columns N/O ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 97.61
$wsTasas.Range("O10").Value = 4140
$wsTasas.Range("N12").Value = 4161.95
$wsTasas.Range("O12").Value = 93.90000000000001
